$d = $word.ActiveDocument

# The document has a few stray empty "List Paragraph" paragraphs (italic,
# sz 28, no visible text) left behind under two of the list items. Remove
# them:
#   - one empty ListParagraph paragraph right after
#     "What are the requirements of the problem?"
#   - two consecutive empty ListParagraph paragraphs (one with no extra
#     indent, one with a 1080-twip left indent) right after
#     "Talk about deployment? Webservice? API? Batch job?"

function Get-ParagraphByText($searchText) {
    $target = $null
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text.TrimEnd("`r", "`a", "`f")
        if ($t -eq $searchText) {
            $target = $p
        }
    }
    return $target
}

# Second location in the document: delete the two stray empty paragraphs,
# bottom one first so the earlier one's identity/range stays valid.
$anchor2 = Get-ParagraphByText("Talk about deployment? Webservice? API? Batch job?")
$stray1 = $anchor2.Next()
$stray2 = $stray1.Next()
$stray2.Range.Delete() | Out-Null
$stray1.Range.Delete() | Out-Null

# First location in the document: delete the single stray empty paragraph.
$anchor1 = Get-ParagraphByText("What are the requirements of the problem?")
$stray3 = $anchor1.Next()
$stray3.Range.Delete() | Out-Null
